# Compare_Tests.xlsx -- add a "Section BLEU Score" template sheet and
# rename the existing sheet to "Whole Article BLEU Score".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the original sheet.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Whole Article BLEU Score"

# ---------------------------------------------------------------------
# 2) Insert the new sheet right after it.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Section BLEU Score"

# ---------------------------------------------------------------------
# 3) Header row for the new sheet -- same template as the Whole Article
#    sheet, but with a "Section Name" column inserted after "Article
#    Name" and the sentence-count headers renamed for sections.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Article Name"
$ws2.Range("B1").Value = "Section Name"
$ws2.Range("C1").Value = "Lang 1"
$ws2.Range("D1").Value = "Lang 2"
$ws2.Range("E1").Value = "Sentences Section 1"
$ws2.Range("F1").Value = "Sentences Section 2"
$ws2.Range("G1").Value = "Est. Avg. Time Per Iteration (sec)"
$ws2.Range("H1").Value = "Num Iterations O(m*n)"
$ws2.Range("I1").Value = "Compare Time Est (min)"
$ws2.Range("J1").Value = "Compare Time Est (sec)"
$ws2.Range("K1").Value = "Compare Time Actual (sec) 1"
$ws2.Range("L1").Value = "Compare Time Actual (sec) 2"
$ws2.Range("M1").Value = "Compare Time Actual (sec) 3"
$ws2.Range("N1").Value = "Compare Time Actual (sec) 4"
$ws2.Range("O1").Value = "Compare Time Actual (sec) 5"
$ws2.Range("P1").Value = "Avg (sec)"
$ws2.Range("Q1").Value = "Avg (min)"
$ws2.Range("R1").Value = "Avg per Iteration (sec)"
$ws2.Range("S1").Value = "Est - Actual"
$ws2.Range("T1").Value = "Open Apps"

# Center-align the "Compare Time Actual / Avg / Est - Actual / Open Apps"
# block, matching the style already used on the first sheet.
$ws2.Range("K1:T1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4) Column widths on the new sheet (character units; the host stores
#    width as ColumnWidth + 5/6, so compensate to land on the desired
#    stored width).
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 16.830729166666668
$ws2.Columns.Item(2).ColumnWidth = 14.998697916666666
$ws2.Columns.Item(5).ColumnWidth = 18.998697916666668
$ws2.Columns.Item(6).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 27.830729166666668
$ws2.Columns.Item(8).ColumnWidth = 19.830729166666668
$ws2.Columns.Item(9).ColumnWidth = 20.666666666666668
$ws2.Columns.Item(10).ColumnWidth = 20.830729166666668
$ws2.Columns.Item(11).ColumnWidth = 29.998697916666668
$ws2.Columns.Item(12).ColumnWidth = 26.666666666666668
$ws2.Columns.Item(13).ColumnWidth = 26.666666666666668
$ws2.Columns.Item(14).ColumnWidth = 25.998697916666668
$ws2.Columns.Item(15).ColumnWidth = 31.998697916666668
$ws2.Columns.Item(18).ColumnWidth = 22.666666666666668
$ws2.Columns.Item(19).ColumnWidth = 11.666666666666666
$ws2.Columns.Item(20).ColumnWidth = 19.998697916666668

# ---------------------------------------------------------------------
# 5) View / selection state.
# ---------------------------------------------------------------------
# Whole Article BLEU Score: selection moves from R14 to P1:S1.
$ws1.Range("P1:S1").Select()

# Section BLEU Score becomes the active (visible) tab, selection S7.
$ws2.Activate()
$ws2.Range("S7").Select()
